$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell such that it is stored as a
# shared-string (t="s") cell even when the text looks numeric, without
# Excel re-interpreting it as a number and without creating a new cell
# style (keeps default style s="0"/absent).
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# --- rows 8-12: make sure every cell in B:F exists (even blank ones) ---
$blankStyleSource = $ws.Range("Z1")
$ws.Range("B8:F12").Style = $blankStyleSource.Style

# Row 8
Set-TextValue $ws.Range("B8") "2"
Set-TextValue $ws.Range("C8") "testX"
Set-TextValue $ws.Range("D8") "3"
Set-TextValue $ws.Range("E8") "4"

# Row 9
Set-TextValue $ws.Range("B9") "5"
Set-TextValue $ws.Range("C9") "testX"
Set-TextValue $ws.Range("D9") "6"
$ws.Range("F9").Value = 7

# Row 10
Set-TextValue $ws.Range("B10") "8"
Set-TextValue $ws.Range("C10") "testX"
Set-TextValue $ws.Range("D10") "9"

# Row 11
Set-TextValue $ws.Range("B11") "10"
Set-TextValue $ws.Range("C11") "testX"
Set-TextValue $ws.Range("E11") "11"

# Row 12
Set-TextValue $ws.Range("B12") "12"
Set-TextValue $ws.Range("C12") "testX"

# Updated result values in rows 4-7 (SimpleRules Collect test table)
Set-TextValue $ws.Range("E7") "23"
Set-TextValue $ws.Range("F6") "24"
Set-TextValue $ws.Range("E6") "25"
Set-TextValue $ws.Range("F4") "27"
Set-TextValue $ws.Range("F7") "2.2"
Set-TextValue $ws.Range("E5") "2.6"

Write-Host "done"
